$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.483.88"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.94%  '
$ws.Range("D3").Value = "'3.462.20"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.99%  '
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").Value = "'583.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +6.04%  '
$ws.Range("D6").Value = "'176.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.29%  '
$ws.Range("E7").Value = '  +3.96%  '
$ws.Range("E8").Value = '  -0.12%  '
$ws.Range("E9").Value = '  +0.46%  '
$ws.Range("E10").Value = '  +5.45%  '
$ws.Range("D11").Value = "'55.48"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.14%  '
$ws.Range("E12").Value = '  +2.85%  '
$ws.Range("E13").Value = '  -1.03%  '
$ws.Range("D14").Value = "'4.025.88"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.27%  '
$ws.Range("D15").Value = "'3.472.87"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.11%  '
$ws.Range("E16").Value = '  +0.20%  '
$ws.Range("D17").Value = "'18.21"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.55%  '
$ws.Range("D18").Value = "'12.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.14%  '
$ws.Range("D19").Value = "'65.496.17"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.11%  '
$ws.Range("E20").Value = '  +1.43%  '
$ws.Range("D21").Value = "'410.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.74%  '
$ws.Range("D22").Value = "'4.29"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +6.26%  '
$ws.Range("D23").Value = "'4.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +7.66%  '
$ws.Range("D24").Value = "'84.44"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.40%  '
$ws.Range("D25").Value = "'13.37"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +10.79%  '
$ws.Range("E26").Value = '  -0.35%  '
$ws.Range("D27").Value = "'2.86"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.87%  '
$ws.Range("D28").Value = "'9.15"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.16%  '
$ws.Range("D29").Value = "'30.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.03%  '
$ws.Range("D30").Value = "'6.73"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.77%  '
$ws.Range("E31").Value = '  -0.05%  '
$ws.Range("D32").Value = "'593.10"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -8.22%  '
$ws.Range("E33").Value = '  -0.88%  '
$ws.Range("D34").Value = "'60.64"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.92%  '
$ws.Range("E35").Value = '  +1.47%  '
$ws.Range("E36").Value = '  +0.08%  '
$ws.Range("E37").Value = '  -3.46%  '
$ws.Range("D38").Value = "'36.78"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.00%  '
$ws.Range("E39").Value = '  +5.97%  '
$ws.Range("E40").Value = '  -1.70%  '
$ws.Range("D41").Value = "'3.218.52"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.17%  '
$ws.Range("E42").Value = '  +0.05%  '
$ws.Range("D43").Value = "'2.96"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.34%  '
$ws.Range("D44").Value = "'3.30"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.53%  '
$ws.Range("E45").Value = '  -4.89%  '
$ws.Range("D46").Value = "'0.0416"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.10%  '
$ws.Range("E47").Value = '  +1.66%  '
$ws.Range("D48").Value = "'2.64"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.46%  '
$ws.Range("D49").Value = "'8.58"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.33%  '
$ws.Range("D50").Value = "'138.47"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.12%  '
$ws.Range("D51").Value = "'2.35"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.24%  '
